$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(7, 6).Value = 83
$ws.Cells.Item(7, 7).Value = 3886.89

$ws.Cells.Item(10, 2).Value = 20465.96

$ws.Cells.Item(37, 6).Value = 157
$ws.Cells.Item(37, 7).Value = 4182.48

$ws.Cells.Item(41, 2).Value = 82498.64

$ws.Cells.Item(49, 6).Value = 107
$ws.Cells.Item(49, 7).Value = 4506.84

$ws.Cells.Item(55, 6).Value = 281
$ws.Cells.Item(55, 7).Value = 15764.1

$ws.Cells.Item(68, 6).Value = 26
$ws.Cells.Item(68, 7).Value = 486.46

$ws.Cells.Item(74, 2).Value = 288349.86

$ws.Cells.Item(115, 6).Value = 23
$ws.Cells.Item(115, 7).Value = 21089.62

$ws.Cells.Item(116, 2).Value = 30798.05

$ws.Cells.Item(124, 6).Value = 86
$ws.Cells.Item(124, 7).Value = 6980.62

$ws.Cells.Item(131, 6).Value = 74
$ws.Cells.Item(131, 7).Value = 4713.8

$ws.Cells.Item(147, 2).Value = 99159.42999999999

$ws.Cells.Item(178, 6).Value = 96
$ws.Cells.Item(178, 7).Value = 4997.76

$ws.Cells.Item(184, 2).Value = 30075.07

$ws.Cells.Item(200, 6).Value = 30
$ws.Cells.Item(200, 7).Value = 2904.3

$ws.Cells.Item(202, 2).Value = -3292.2

$ws.Cells.Item(224, 6).Value = 100
$ws.Cells.Item(224, 7).Value = 7861

$ws.Cells.Item(234, 6).Value = 67
$ws.Cells.Item(234, 7).Value = 5999.18

$ws.Cells.Item(236, 2).Value = 63255
$ws.Cells.Item(236, 6).Value = 93
$ws.Cells.Item(236, 7).Value = 7626

$ws.Cells.Item(237, 2).Value = 57004
$ws.Cells.Item(237, 6).Value = 5
$ws.Cells.Item(237, 7).Value = 410

$ws.Cells.Item(241, 2).Value = 64329
$ws.Cells.Item(241, 5).Value = 128.32
$ws.Cells.Item(241, 6).Value = 1
$ws.Cells.Item(241, 7).Value = 120.69

$ws.Cells.Item(242, 2).Value = 57552
$ws.Cells.Item(242, 5).Value = 136.86
$ws.Cells.Item(242, 6).Value = -5
$ws.Cells.Item(242, 7).Value = -603.45

$ws.Cells.Item(248, 6).Value = 42
$ws.Cells.Item(248, 7).Value = 2515.38

$ws.Cells.Item(250, 2).Value = 101796.41

$ws.Cells.Item(275, 6).Value = 169
$ws.Cells.Item(275, 7).Value = 9195.290000000001

$ws.Cells.Item(283, 2).Value = 115610.45

$ws.Cells.Item(287, 6).Value = 2281
$ws.Cells.Item(287, 7).Value = 42198.5

$ws.Cells.Item(294, 2).Value = 64421.6

$ws.Cells.Item(325, 2).Value = 66188
$ws.Cells.Item(325, 3).Value = "HIM-Baby Care Gift Pack (Ww)1"
$ws.Cells.Item(325, 4).Value = 315.8
$ws.Cells.Item(325, 5).Value = 377.31
$ws.Cells.Item(325, 6).Value = 49
$ws.Cells.Item(325, 7).Value = 15474.2

$ws.Cells.Item(326, 2).Value = 48719
$ws.Cells.Item(326, 3).Value = "HIM-BABY CARE GIFT PACK (WW)1"
$ws.Cells.Item(326, 4).Value = 295.75
$ws.Cells.Item(326, 5).Value = 353.35
$ws.Cells.Item(326, 6).Value = -82
$ws.Cells.Item(326, 7).Value = -24251.5

$ws.Cells.Item(362, 6).Value = 229
$ws.Cells.Item(362, 7).Value = 10733.23

$ws.Cells.Item(363, 6).Value = 26
$ws.Cells.Item(363, 7).Value = 2884.44

$ws.Cells.Item(366, 6).Value = 124
$ws.Cells.Item(366, 7).Value = 6787.76

$ws.Cells.Item(375, 2).Value = 179053.58

$ws.Cells.Item(400, 2).Value = 60325
$ws.Cells.Item(400, 5).Value = 151.57
$ws.Cells.Item(400, 6).Value = -102
$ws.Cells.Item(400, 7).Value = -12939.72

$ws.Cells.Item(401, 2).Value = 63560
$ws.Cells.Item(401, 5).Value = 134.87
$ws.Cells.Item(401, 6).Value = 1
$ws.Cells.Item(401, 7).Value = 126.86

$ws.Cells.Item(404, 6).Value = 42
$ws.Cells.Item(404, 7).Value = 4364.22

$ws.Cells.Item(405, 6).Value = 39
$ws.Cells.Item(405, 7).Value = 5374.98

$ws.Cells.Item(408, 6).Value = 213
$ws.Cells.Item(408, 7).Value = 36493.29

$ws.Cells.Item(412, 2).Value = 52249.63

$ws.Cells.Item(442, 6).Value = 2
$ws.Cells.Item(442, 7).Value = 86.02

$ws.Cells.Item(454, 2).Value = 99279.39

$ws.Cells.Item(471, 6).Value = 343
$ws.Cells.Item(471, 7).Value = 56934.57

$ws.Cells.Item(473, 2).Value = 100955.83

$ws.Cells.Item(483, 2).Value = 58047
$ws.Cells.Item(483, 4).Value = 105.54
$ws.Cells.Item(483, 5).Value = 126.1
$ws.Cells.Item(483, 6).Value = 34
$ws.Cells.Item(483, 7).Value = 3588.36

$ws.Cells.Item(484, 2).Value = 47097
$ws.Cells.Item(484, 4).Value = 112.28
$ws.Cells.Item(484, 5).Value = 134.16
$ws.Cells.Item(484, 6).Value = 15
$ws.Cells.Item(484, 7).Value = 1684.2

$ws.Cells.Item(489, 6).Value = 642
$ws.Cells.Item(489, 7).Value = 62017.2

$ws.Cells.Item(492, 2).Value = 80840.39999999999

$ws.Cells.Item(508, 6).Value = 22
$ws.Cells.Item(508, 7).Value = 521.4

$ws.Cells.Item(517, 6).Value = 149
$ws.Cells.Item(517, 7).Value = 8577.93

$ws.Cells.Item(520, 6).Value = 80
$ws.Cells.Item(520, 7).Value = 4376

$ws.Cells.Item(522, 2).Value = 207262.6

$ws.Cells.Item(525, 6).Value = 40
$ws.Cells.Item(525, 7).Value = 6356.8

$ws.Cells.Item(534, 2).Value = 25810.5

$ws.Cells.Item(567, 2).Value = 64925
$ws.Cells.Item(567, 5).Value = 13.97
$ws.Cells.Item(567, 6).Value = 111
$ws.Cells.Item(567, 7).Value = 1459.65

$ws.Cells.Item(568, 2).Value = 45709
$ws.Cells.Item(568, 5).Value = 15.69
$ws.Cells.Item(568, 6).Value = -300
$ws.Cells.Item(568, 7).Value = -3945

$ws.Cells.Item(611, 6).Value = 142
$ws.Cells.Item(611, 7).Value = 40180.32

$ws.Cells.Item(615, 2).Value = 150026.05

$ws.Cells.Item(653, 6).Value = 18
$ws.Cells.Item(653, 7).Value = 614.7

$ws.Cells.Item(659, 2).Value = 5447.7

$ws.Cells.Item(662, 2).Value = 64833
$ws.Cells.Item(662, 5).Value = 34.9
$ws.Cells.Item(662, 6).Value = 90
$ws.Cells.Item(662, 7).Value = 2954.7

$ws.Cells.Item(663, 2).Value = 60025
$ws.Cells.Item(663, 5).Value = 37.22
$ws.Cells.Item(663, 6).Value = -98
$ws.Cells.Item(663, 7).Value = -3217.34

$ws.Cells.Item(680, 6).Value = 378
$ws.Cells.Item(680, 7).Value = 37750.86

$ws.Cells.Item(685, 6).Value = 27
$ws.Cells.Item(685, 7).Value = 2848.77

$ws.Cells.Item(688, 6).Value = 559
$ws.Cells.Item(688, 7).Value = 30677.92

$ws.Cells.Item(692, 6).Value = 212
$ws.Cells.Item(692, 7).Value = 18149.32

$ws.Cells.Item(694, 6).Value = 44
$ws.Cells.Item(694, 7).Value = 3886.96

$ws.Cells.Item(695, 2).Value = 187654.97

$ws.Cells.Item(741, 6).Value = 11
$ws.Cells.Item(741, 7).Value = 911.79

$ws.Cells.Item(742, 2).Value = 50591.65

$ws.Cells.Item(806, 6).Value = 1
$ws.Cells.Item(806, 7).Value = 108.81

$ws.Cells.Item(807, 6).Value = 157
$ws.Cells.Item(807, 7).Value = 17083.17

$ws.Cells.Item(812, 6).Value = 42
$ws.Cells.Item(812, 7).Value = 6145.86

$ws.Cells.Item(815, 6).Value = 104
$ws.Cells.Item(815, 7).Value = 15088.32

$ws.Cells.Item(819, 6).Value = 105
$ws.Cells.Item(819, 7).Value = 5052.6

$ws.Cells.Item(823, 6).Value = 22
$ws.Cells.Item(823, 7).Value = 727.3200000000001

$ws.Cells.Item(828, 6).Value = 106
$ws.Cells.Item(828, 7).Value = 54541.24

$ws.Cells.Item(831, 6).Value = 430
$ws.Cells.Item(831, 7).Value = 15836.9

$ws.Cells.Item(838, 2).Value = 333786.76

$ws.Cells.Item(874, 6).Value = 91
$ws.Cells.Item(874, 7).Value = 4872.14

$ws.Cells.Item(878, 6).Value = 101
$ws.Cells.Item(878, 7).Value = 8111.31

$ws.Cells.Item(879, 6).Value = 79
$ws.Cells.Item(879, 7).Value = 2825.83

$ws.Cells.Item(885, 2).Value = 25778.41

$ws.Cells.Item(890, 6).Value = 244
$ws.Cells.Item(890, 7).Value = 7376.12

$ws.Cells.Item(897, 2).Value = 347202.72

$ws.Cells.Item(915, 6).Value = 14
$ws.Cells.Item(915, 7).Value = 3393.46

$ws.Cells.Item(931, 6).Value = 141
$ws.Cells.Item(931, 7).Value = 5200.08

$ws.Cells.Item(933, 6).Value = 23
$ws.Cells.Item(933, 7).Value = 4923.38

$ws.Cells.Item(936, 2).Value = 118723.16

$ws.Cells.Item(942, 2).Value = 5147886.35

$ws.Cells.Item(943, 2).Value = 5147886.35
